$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost (Easting) / Nord (Northing) coordinates to whole numbers
$ws.Range("Q2").Value = 572178
$ws.Range("R2").Value = 6708775

# Clear the Starttid / Sluttid cells (time-of-day values removed)
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
